# Add a blank placeholder data row (row 2) below the header row.
# Text / lookup columns are left blank; amount columns default to 0,
# matching the report generator's "no records" row shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the whole row first so every column gets an explicit (blank) value.
$ws.Range("A2:T2").Value = ""

# Amount / numeric columns default to 0:
#   I = Đơn giá gốc, K = Upsale, L = Đơn giá, M = Thanh toán lần đầu,
#   N = Trả sau, O = Đã thanh toán, P = Dư nợ
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0

# B = Mã dịch vụ is a numeric column too, but stays blank (0) for this row.
$ws.Range("B2").Value = 0
